# daily auto push: 2026-01-07 02:26 UTC
#
# The sheet logs forecast rows of (date, weekday, hour, rank). The most
# recent tracked day is 2026/01/07 (水), whose last logged hour is 3:00 in
# row 570. This push appends two more hourly samples for that same day
# (07:00 and 10:00), which get inserted directly above the older
# "2026/12/29" block, pushing every row from the old 571 onward down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new samples right after the existing 2026/01/07 rows.
$ws.Rows("571:572").Insert()

# Seed the date/weekday text by copying it (as values) from the row above
# instead of typing the literal "2026/01/07" string - that keeps the cells
# stored as plain text (matching every other row in the column) rather than
# letting Excel auto-convert the date-shaped text into a real date serial.
$ws.Range("A570:B570").Copy()
$ws.Range("A571:B571").PasteSpecial(-4163)
$ws.Range("A570:B570").Copy()
$ws.Range("A572:B572").PasteSpecial(-4163)

# New hourly readings for 2026/01/07 (水).
$ws.Cells.Item(571, 3).Value = 7
$ws.Cells.Item(571, 4).Value = 20

$ws.Cells.Item(572, 3).Value = 10
$ws.Cells.Item(572, 4).Value = 20
